# Update NB3 Body PCB (rev 6)
# Updates the BOM on the NB3_power sheet:
#  - Row2 (J12,J13): switch 2x03 header part from SMD (Ckmtw) to a
#    through-hole part from XFCN, updates mfr#, supplier# and unit price.
#  - Row3 (2x05 headers): mfr# column now carries the "2.54-2*5P" part code
#    instead of repeating the supplier number.
#  - Row4 (2x07 headers): mfr# column now carries the "2.54-2*7P" part code
#    instead of repeating the supplier number.
#  - Row5 (PH header): designator list grows to include J11, so quantity
#    goes from 1 to 2; matches the font/alignment used on row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: J12,J13 2x03 header swapped to through-hole part ---
$ws.Range("C2").Value = "Male Header 2x03 Through-Hole 2.54mm P=2.54mm"
$ws.Range("D2").Value = "XFCN"
$ws.Range("E2").Value = "PZ254V-12-6P"
$ws.Range("G2").Value = "C492420"
$ws.Range("H2").Value = 0.0354

# --- Row 3: 2x05 header, mfr number column updated ---
$ws.Range("E3").Value = "2.54-2*5P"

# --- Row 4: 2x07 header, mfr number column updated ---
$ws.Range("E4").Value = "2.54-2*7P"

# --- Row 5: PH header now also used on J11 ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "J11,J10"
$ws.Range("B5").Font.Name = "Noto Sans"
$ws.Range("B5").HorizontalAlignment = 1

# --- View state: selection moved to G10 ---
$ws.Range("G10").Select() | Out-Null
